$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 5) for a 4th shareholder. Its identity / contact
# columns (B, C, F, G, H) mirror the first shareholder's row exactly, so
# copy-paste the values from row 2 (this reuses the existing shared
# strings and avoids Excel re-interpreting the leading-zero numeric-looking
# text, e.g. "001090001234" or "0912345678", as a number).
$ws.Range("A2:H2").Copy()
$ws.Range("A5:H5").PasteSpecial(-4163)  # xlPasteValues

# Now set the columns that differ for this new shareholder: running
# number, shareholder code (a brand-new code, SH004) and share count.
$ws.Range("A5").Value = 4
$ws.Range("D5").Value = "SH004"
$ws.Range("E5").Value = 200

# I5 holds an issue date (serial 42040 => 2015-02-05). Reuse the date
# style already applied to I2:I4 (copy format from I4) and then write the
# raw date serial so no new number format / style gets created.
$ws.Range("I4").Copy()
$ws.Range("I5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I5").Value = 42040

# Update the active selection, matching the post-edit workbook state.
$ws.Range("K9").Select()
